# Excel COM-interop script implementing the commit:
#   opt@brand#:add permission control config.
#
# Sheet 1 ("需要新增加的权限节点（brand）") gets a new permission-config
# data row (row 14), copied down from the existing data row above it so it
# keeps the same cell formatting, then filled in with the new permission
# details. The previously-selected cell (C19) also shifts up to C18 to
# match where the author's selection ended up when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-assert row 13's permission-path label (B13) as part of finishing the
# new block of rows.
$ws.Range("B13").Value = "品牌管理-优先品牌管理-品牌评测管理"

# Bring row 14 into the same formatting as the other populated data rows
# (row 11 uses the "clean" s=3 / s=14 pattern that the new row should use).
$ws.Range("A11:I11").Copy()
$ws.Range("A14:I14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14").ClearFormats()           # A14 keeps the default (unstyled) look

# Fill in the new permission-config row.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "品牌管理-优先品牌管理-品牌评测管理"
$ws.Range("C14").Value = "cms/addContextForUploadPic"
$ws.Range("D14").Value = "编辑器上传文件"
$ws.Range("E14").Value = "否"
$ws.Range("F14").Value = "是"
$ws.Range("G14").Value = "bs"
$ws.Range("H14").Value = "EQ组"

# Move the active selection to C18, matching the saved view state.
$ws.Range("C18").Select()
